$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("X2").Value = 0.1806907185664647
$ws.Range("Y2").Value = -0.09356977516318091
$ws.Range("AB2").Value = 0.08384283964516498
$ws.Range("AC2").Value = -0.08384283964516498
$ws.Range("X3").Value = 0.06821556785261536
$ws.Range("Y3").Value = 0.05996488970564849
$ws.Range("AB3").Value = 0.06858016103713035
$ws.Range("AC3").Value = -0.06858016103713035
$ws.Range("X4").Value = 0.128514364377689
$ws.Range("Y4").Value = -0.265551401414726
$ws.Range("AB4").Value = 0.07906163593284338
$ws.Range("AC4").Value = -0.07906163593284338
$ws.Range("X5").Value = 0.125613453732096
$ws.Range("Y5").Value = -0.06606879523822558
$ws.Range("AB5").Value = 0.08206669324416022
$ws.Range("AC5").Value = -0.08206669324416022
$ws.Range("X6").Value = 0.1704897576835411
$ws.Range("Y6").Value = -0.08903832159496401
$ws.Range("AB6").Value = 0.08224455118508134
$ws.Range("AC6").Value = -0.08224455118508134
$ws.Range("X7").Value = 0.1307826060119376
$ws.Range("Y7").Value = -0.03982374390106132
$ws.Range("AB7").Value = 0.08269891802136012
$ws.Range("AC7").Value = -0.08269891802136012
$ws.Range("X8").Value = 0.2456543363469879
$ws.Range("Y8").Value = -0.1623713116512966
$ws.Range("AB8").Value = 0.08498676126896983
$ws.Range("AC8").Value = -0.08498676126896983
$ws.Range("X9").Value = 0.4263125493242699
$ws.Range("Y9").Value = -0.2656513443652327
$ws.Range("AB9").Value = 0.08742894766626919
$ws.Range("AC9").Value = -0.08742894766626919
$ws.Range("X10").Value = 0.4605873926857427
$ws.Range("Y10").Value = -0.361699739738023
$ws.Range("AB10").Value = 0.08766630616089224
$ws.Range("AC10").Value = -0.08766630616089224
$ws.Range("X11").Value = 0.2651739770449885
$ws.Range("Y11").Value = -0.1656621560646918
$ws.Range("AB11").Value = 0.1056358876986679
$ws.Range("AC11").Value = -0.1056358876986679
$ws.Range("X12").Value = 0.1908916794493883
$ws.Range("Y12").Value = -0.1413747712368279
$ws.Range("AB12").Value = 0.1067806737634219
$ws.Range("AC12").Value = -0.1067806737634219
